$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 215.70833
$ws.Range("I33").Value = 95.61539
$ws.Range("J33").Value = 357.63635
$ws.Range("K33").Value = 95.61539
$ws.Range("L33").Value = 357.63635
$ws.Range("M33").Value = 133.38461
$ws.Range("N33").Value = -815.63635
# Row 64
$ws.Range("H64").Value = 4697.778
$ws.Range("I64").Value = 4546.6665
$ws.Range("K64").Value = 4546.6665
$ws.Range("M64").Value = -4298.6665
# Row 67
$ws.Range("H67").Value = 4697.778
$ws.Range("I67").Value = 4546.6665
$ws.Range("K67").Value = 4546.6665
$ws.Range("M67").Value = -3688.6665
# Row 99
$ws.Range("H99").Value = 200
$ws.Range("I99").Value = 200
$ws.Range("K99").Value = 600
$ws.Range("M99").Value = 898
# Row 137
$ws.Range("H137").Value = 847.2558
$ws.Range("I137").Value = 790
$ws.Range("J137").Value = 995.1667
$ws.Range("K137").Value = 2370
$ws.Range("L137").Value = 2985.5001
$ws.Range("M137").Value = 180
$ws.Range("N137").Value = -8085.5001
# Row 138
$ws.Range("H138").Value = 3529.4062
$ws.Range("I138").Value = 1797.1
$ws.Range("J138").Value = 4316.8184
$ws.Range("K138").Value = 5391.299999999999
$ws.Range("L138").Value = 12950.4552
$ws.Range("M138").Value = -251.2999999999993
$ws.Range("N138").Value = -23230.4552
# Row 141
$ws.Range("H141").Value = 797.5
$ws.Range("I141").Value = 797.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2392.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2787.5
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
# Row 18
$ws.Range("H18").Value = 100013
$ws.Range("J18").Value = 100013
$ws.Range("L18").Value = 100013
$ws.Range("N18").Value = -100657
# Row 32
$ws.Range("H32").Value = 4090.79
$ws.Range("I32").Value = 4062.6702
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 4062.6702
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -3775.6702
$ws.Range("N32").Value = -5574
# Row 88
$ws.Range("H88").Value = 2673.3125
$ws.Range("I88").Value = 2044.875
$ws.Range("J88").Value = 3301.75
$ws.Range("K88").Value = 2044.875
$ws.Range("L88").Value = 3301.75
$ws.Range("M88").Value = -1638.875
$ws.Range("N88").Value = -4113.75
# Row 91
$ws.Range("H91").Value = 2673.3125
$ws.Range("I91").Value = 2044.875
$ws.Range("J91").Value = 3301.75
$ws.Range("K91").Value = 2044.875
$ws.Range("L91").Value = 3301.75
$ws.Range("M91").Value = -640.875
$ws.Range("N91").Value = -6109.75
# Row 102
$ws.Range("H102").Value = 5567.615
$ws.Range("I102").Value = 3439.9
$ws.Range("K102").Value = 3439.9
$ws.Range("M102").Value = -1817.9
# Row 132
$ws.Range("H132").Value = 1439.9531
$ws.Range("I132").Value = 1168.4131
$ws.Range("J132").Value = 2133.889
$ws.Range("K132").Value = 3505.2393
$ws.Range("L132").Value = 6401.667
$ws.Range("M132").Value = -975.2393000000002
$ws.Range("N132").Value = -11461.667

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 29126.541
$ws.Range("I134").Value = 1915.5483
$ws.Range("K134").Value = 5746.644899999999
$ws.Range("M134").Value = -3211.644899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 12000000
$ws.Range("I12").Value = 12000000
$ws.Range("K12").Value = 12000000
$ws.Range("M12").Value = -11999830
# Row 16
$ws.Range("H16").Value = 1704.7222
$ws.Range("I16").Value = 1437.2
$ws.Range("J16").Value = 2039.125
$ws.Range("K16").Value = 1437.2
$ws.Range("L16").Value = 2039.125
$ws.Range("M16").Value = -1150.2
$ws.Range("N16").Value = -2613.125
# Row 31
$ws.Range("H31").Value = 32385.912
$ws.Range("I31").Value = 3359.3044
$ws.Range("J31").Value = 93077.91
$ws.Range("K31").Value = 3359.3044
$ws.Range("L31").Value = 93077.91
$ws.Range("M31").Value = -3064.3044
$ws.Range("N31").Value = -93667.91
# Row 34
$ws.Range("H34").Value = 32385.912
$ws.Range("I34").Value = 3359.3044
$ws.Range("J34").Value = 93077.91
$ws.Range("K34").Value = 3359.3044
$ws.Range("L34").Value = 93077.91
$ws.Range("M34").Value = -3157.3044
$ws.Range("N34").Value = -93481.91
# Row 44
$ws.Range("H44").Value = 20450
$ws.Range("J44").Value = 20450
$ws.Range("L44").Value = 20450
$ws.Range("N44").Value = -21334
# Row 105
$ws.Range("H105").Value = 6830
$ws.Range("I105").Value = 7996
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 7996
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = -6249
$ws.Range("N105").Value = -4494
# Row 113
$ws.Range("H113").Value = 1704.7222
$ws.Range("I113").Value = 1437.2
$ws.Range("J113").Value = 2039.125
$ws.Range("K113").Value = 1437.2
$ws.Range("L113").Value = 2039.125
$ws.Range("M113").Value = 732.8
$ws.Range("N113").Value = -6379.125
# Row 132
$ws.Range("H132").Value = 2603.25
$ws.Range("I132").Value = 4999
$ws.Range("J132").Value = 1804.6666
$ws.Range("K132").Value = 14997
$ws.Range("L132").Value = 5413.9998
$ws.Range("M132").Value = -12467
$ws.Range("N132").Value = -10473.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 3730
$ws.Range("J32").Value = 3730
$ws.Range("L32").Value = 11190
$ws.Range("N32").Value = -11756
# Row 80
$ws.Range("H80").Value = 8385.929
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 8385.929
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 25157.787
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -27029.787
# Row 83
$ws.Range("H83").Value = 8385.929
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 8385.929
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 75473.361
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -84833.361
# Row 121
$ws.Range("H121").Value = 824.2778
$ws.Range("J121").Value = 852.17645
$ws.Range("L121").Value = 2556.52935
$ws.Range("N121").Value = -5176.529350000001

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 3997.5
$ws.Range("I43").Value = 3997.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3997.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3846.5
$ws.Range("N43").ClearContents()
# Row 46
$ws.Range("H46").Value = 11273.75
$ws.Range("I46").Value = 595
$ws.Range("J46").Value = 14833.333
$ws.Range("K46").Value = 595
$ws.Range("L46").Value = 14833.333
$ws.Range("M46").Value = -439
$ws.Range("N46").Value = -15145.333
# Row 70
$ws.Range("H70").Value = 32598.285
$ws.Range("I70").Value = 53296
$ws.Range("J70").Value = 5001.3335
$ws.Range("K70").Value = 53296
$ws.Range("L70").Value = 5001.3335
$ws.Range("M70").Value = -53026
$ws.Range("N70").Value = -5541.3335
# Row 73
$ws.Range("H73").Value = 32598.285
$ws.Range("I73").Value = 53296
$ws.Range("J73").Value = 5001.3335
$ws.Range("K73").Value = 53296
$ws.Range("L73").Value = 5001.3335
$ws.Range("M73").Value = -52360
$ws.Range("N73").Value = -6873.3335
# Row 132
$ws.Range("H132").Value = 3316.963
$ws.Range("I132").Value = 3316.4614
$ws.Range("J132").Value = 3317.4285
$ws.Range("K132").Value = 9949.3842
$ws.Range("L132").Value = 9952.2855
$ws.Range("M132").Value = -7419.3842
$ws.Range("N132").Value = -15012.2855
# Row 134
$ws.Range("H134").Value = 16662
$ws.Range("J134").Value = 16662
$ws.Range("L134").Value = 49986
$ws.Range("N134").Value = -55056
# Row 135
$ws.Range("H135").Value = 58886.668
$ws.Range("J135").Value = 58886.668
$ws.Range("L135").Value = 58886.668
$ws.Range("N135").Value = -69026.66800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 111
$ws.Range("H111").Value = 36693.5
$ws.Range("J111").Value = 36693.5
$ws.Range("L111").Value = 36693.5
$ws.Range("N111").Value = -44873.5
# Row 132
$ws.Range("H132").Value = 3358.7646
$ws.Range("I132").Value = 2200
$ws.Range("J132").Value = 4273.579
$ws.Range("K132").Value = 6600
$ws.Range("L132").Value = 12820.737
$ws.Range("M132").Value = -4070
$ws.Range("N132").Value = -17880.737

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 976.0862
$ws.Range("I132").Value = 813.14636
$ws.Range("J132").Value = 1369.0588
$ws.Range("K132").Value = 2439.43908
$ws.Range("L132").Value = 4107.1764
$ws.Range("M132").Value = 90.5609199999999
$ws.Range("N132").Value = -9167.1764
# Row 136
$ws.Range("H136").Value = 2523.6
$ws.Range("I136").Value = 2787.8333
$ws.Range("J136").Value = 1466.6666
$ws.Range("K136").Value = 8363.499899999999
$ws.Range("L136").Value = 4399.9998
$ws.Range("M136").Value = -5813.499899999999
$ws.Range("N136").Value = -9499.9998
